{"js": "// Convert the five Russian use-case heading labels to uppercase, leaving\n// everything else (formatting, other runs, paragraphs) untouched.\nconst replacements = [\n  [\"\u041a\u0440\u0430\u0442\u043a\u043e\u0435 \u043e\u043f\u0438\u0441\u0430\u043d\u0438\u0435\", \"\u041a\u0420\u0410\u0422\u041a\u041e\u0415 \u041e\u041f\u0418\u0421\u0410\u041d\u0418\u0415\"],\n  [\"\u041e\u0441\u043d\u043e\u0432\u043d\u043e\u0439 \u043f\u043e\u0442\u043e\u043a \u0441\u043e\u0431\u044b\u0442\u0438\u0439\", \"\u041e\u0421\u041d\u041e\u0412\u041d\u041e\u0419 \u041f\u041e\u0422\u041e\u041a \u0421\u041e\u0411\u042b\u0422\u0418\u0419\"],\n  [\"\u0410\u043b\u044c\u0442\u0435\u0440\u043d\u0430\u0442\u0438\u0432\u043d\u044b\u0435 \u043f\u043e\u0442\u043e\u043a\u0438\", \"\u0410\u041b\u042c\u0422\u0415\u0420\u041d\u0410\u0422\u0418\u0412\u041d\u042b\u0415 \u041f\u041e\u0422\u041e\u041a\u0418\"],\n  [\"\u041f\u0440\u0435\u0434\u0443\u0441\u043b\u043e\u0432\u0438\u044f\", \"\u041f\u0420\u0415\u0414\u0423\u0421\u041b\u041e\u0412\u0418\u042f\"],\n  [\"\u041f\u043e\u0441\u0442\u0443\u0441\u043b\u043e\u0432\u0438\u044f\", \"\u041f\u041e\u0421\u0422\u0423\u0421\u041b\u041e\u0412\u0418\u042f\"],\n];\n\nfor (const [findText, newText] of replacements) {\n  const results = context.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items/text\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Convert the five Russian use-case heading labels to uppercase, leaving\n# everything else (formatting, other runs, paragraphs) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"\u041a\u0440\u0430\u0442\u043a\u043e\u0435 \u043e\u043f\u0438\u0441\u0430\u043d\u0438\u0435\"; Replace = \"\u041a\u0420\u0410\u0422\u041a\u041e\u0415 \u041e\u041f\u0418\u0421\u0410\u041d\u0418\u0415\" },\n    @{ Find = \"\u041e\u0441\u043d\u043e\u0432\u043d\u043e\u0439 \u043f\u043e\u0442\u043e\u043a \u0441\u043e\u0431\u044b\u0442\u0438\u0439\"; Replace = \"\u041e\u0421\u041d\u041e\u0412\u041d\u041e\u0419 \u041f\u041e\u0422\u041e\u041a \u0421\u041e\u0411\u042b\u0422\u0418\u0419\" },\n    @{ Find = \"\u0410\u043b\u044c\u0442\u0435\u0440\u043d\u0430\u0442\u0438\u0432\u043d\u044b\u0435 \u043f\u043e\u0442\u043e\u043a\u0438\"; Replace = \"\u0410\u041b\u042c\u0422\u0415\u0420\u041d\u0410\u0422\u0418\u0412\u041d\u042b\u0415 \u041f\u041e\u0422\u041e\u041a\u0418\" },\n    @{ Find = \"\u041f\u0440\u0435\u0434\u0443\u0441\u043b\u043e\u0432\u0438\u044f\"; Replace = \"\u041f\u0420\u0415\u0414\u0423\u0421\u041b\u041e\u0412\u0418\u042f\" },\n    @{ Find = \"\u041f\u043e\u0441\u0442\u0443\u0441\u043b\u043e\u0432\u0438\u044f\"; Replace = \"\u041f\u041e\u0421\u0422\u0423\u0421\u041b\u041e\u0412\u0418\u042f\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($r.Find, $true, $true, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
